# Finished Week 13 logging
# Update row 3 ("R" - Road) target depth stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 247
$wsOff.Range("C3").Value = 164
$wsOff.Range("D3").Value = 85
$wsOff.Range("E3").Value = 38

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 470
$wsDef.Range("C3").Value = 347
$wsDef.Range("D3").Value = 73
$wsDef.Range("E3").Value = 36
$wsDef.Range("G3").Value = 2
